# Automatic update of files.
# Every data row's "Förändrad" (column C) date advances by one day:
# serial 45203 (2023-10-04) -> 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlUp = -4162
$xlUp = -4162

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End($xlUp).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

$oldValue = 45203
$newValue = 45204

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
